$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple price/volume updates (rows with same coin) ---
# D-column values must stay as literal text (preserve trailing zeros, thousand-dot
# notation, etc.), so force text format before assigning, then clear the residual
# number-format style to match the original (unstyled) cells.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.350.43"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.785.47"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "612.30"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "187.20"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +14.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.770.97"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.78%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.642"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.76%  "
$ws.Range("E9").Value = "  +0.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.737"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("E11").Value = "  -3.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "57.63"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +10.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000302"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.74"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.377.59"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.779.87"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.74"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -4.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.15"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.15"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.55%  "
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.022.36"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "421.75"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.52%  "
$ws.Range("E23").Value = "  +1.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "90.80"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.10"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.12"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.33"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.08"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.10"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.73"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.52"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.57"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.71"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.121"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "44.93"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "621.23"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "65.30"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0928"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.412"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.51%  "
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.140"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.84"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.12"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.50"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.829.72"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.25%  "


# --- Row content swaps (coin rows reordered) ---
# Rows 45 & 46 swap: VeChain <-> dogwifhat
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.07"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.32%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0451"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.35%  "

# Rows 48, 50, 51 rotate: Stellar -> ApeXProtocol(48); WEMIXToken -> Stellar(50); ApeXProtocol -> WEMIXToken(51)
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.28"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.93%  "

$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.137"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.71%  "

$ws.Range("B51").Value = "WEMIXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.76"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.02%  "
